# Remove four account rows from the "Export" sheet (Saldo.xlsx):
#   003641655 MARCELO    50000      (row 2)
#   004948033 GUILHERME  30002.74   (row 3)
#   004920447 MARILIA    5000       (row 9, before the earlier deletions shift it)
#   004381180 HFR        4900       (row 10, before the earlier deletions shift it)
#
# Delete from the bottom up so earlier row numbers stay valid as rows shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(10).Delete()
$ws.Rows(9).Delete()
$ws.Rows(3).Delete()
$ws.Rows(2).Delete()
